# Weekly price-table update: a new week's price record is inserted as a
# new row right before the existing row 136 ("Primera" quality, same
# market/category metadata), and every row that used to follow shifts
# down by one (dimension grows from A1:R205 to A1:R206).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new blank row at row 136 - pushes old rows 136..205 down to 137..206
$ws.Rows.Item(136).Insert()

# Seed the new row with the same "metadata" columns (market, region, codes,
# category, variety, quality, unit, origin, Kg/unit flag, classification)
# as the row immediately below it (the row that used to be 136), since the
# new data point belongs to the same series.
$srcRow = $ws.Range("A137:R137")
$dstRow = $ws.Range("A136:R136")
$dstRow.Value2 = $srcRow.Value2

# Now overwrite the columns that actually hold this week's new figures:
# Fecha (date serial), Volumen, Precio minimo, Precio maximo,
# Precio promedio ponderado, Precio $/Kg.
$ws.Range("D136").Value2 = 45016
$ws.Range("J136").Value2 = 150
$ws.Range("K136").Value2 = 1500
$ws.Range("L136").Value2 = 1500
$ws.Range("M136").Value2 = 1500
$ws.Range("P136").Value2 = 1500
